$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell while preserving the cell's
# existing style (Range.Value auto-converts numeric-looking strings like
# "12" into numbers, which loses the shared-string/text type and drags in a
# brand new "@" text style; going through a scratch cell + formula + copy /
# PasteSpecial(values) keeps the original style id and stores the value as
# text, matching how the sheet already stores every figure as text).
function Set-TextValue {
    param($addr, [string]$text)
    $escaped = $text.Replace('"', '""')
    $scratch = $ws.Range("Z200")
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# Week range label
Set-TextValue "B5" "7/2/2018-7/8/2018"

# Day-of-week column headers
Set-TextValue "E7" "2-Jul"
Set-TextValue "F7" "3-Jul"
Set-TextValue "G7" "4-Jul"
Set-TextValue "H7" "5-Jul"
Set-TextValue "I7" "6-Jul"
Set-TextValue "J7" "7-Jul"
Set-TextValue "K7" "8-Jul"

# First activity row: renamed and hours unchanged (still 1/day)
Set-TextValue "B9" "Testing ( 1 ) "

# Second activity row: newly filled in
Set-TextValue "B10" "Client call ( 12 ) "
Set-TextValue "E10" "12"
Set-TextValue "F10" "1"
Set-TextValue "G10" "11"
Set-TextValue "H10" "1"
Set-TextValue "I10" "11"
Set-TextValue "J10" "1"
Set-TextValue "K10" "1"

# Daily totals row, recomputed from the two activity rows above
Set-TextValue "E18" "13"
Set-TextValue "F18" "2"
Set-TextValue "G18" "12"
Set-TextValue "H18" "2"
Set-TextValue "I18" "12"
Set-TextValue "J18" "2"
Set-TextValue "K18" "2"

# Total hours for the week
Set-TextValue "D19" "45"

# Clear the scratch cell so it leaves no trace in the saved workbook
$ws.Range("Z200").ClearContents()

# PasteSpecial into the anchor cell of a merged range drops the merge in
# this runtime, so restore the merges that covered the cells we just wrote.
$ws.Range("B5:D5").Merge()
$ws.Range("B9:C9").Merge()
$ws.Range("B10:C10").Merge()
$ws.Range("D19:L19").Merge()
